$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.364.92"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.036.20"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "384.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "3.525.22"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "3.045.43"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.967"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("E18").Value = "  -5.69%  "
$ws.Range("D19").Value = "51.439.50"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("E26").Value = "  +5.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.170"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.288"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +4.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "2.024.73"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "3.342.70"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.69%  "
